# "quizhou infested city update" -- log update edits for
# "logs and suggested ways/university party.docx"
#
# This script performs the substantive content/structure changes:
#   1. "Going in tower..." -> "Ventuum tower complete", and the trailing
#      _GoBack bookmark that used to sit on that paragraph is relocated
#      (see step 2); a new, empty bulleted list paragraph is appended
#      right after it.
#   2. The _GoBack bookmark is (re)placed between "session (" and the
#      closing ")" a few paragraphs down.
#   3. "...Last redacted 13.05.19)" -> ")" (the "Last redacted 13.05.19"
#      text is dropped, leaving just the closing paren).
#   4. "Disrupt new church building ... (lasts 3 months)" is blanked out
#      to a single space.
#   5. The three quest-line bullets that followed it ("Zingar
#      questline...", "Help Alliance's allies...", "Help halfling
#      allies ") are removed entirely.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Relocate the _GoBack bookmark from the "Going in tower..." paragraph
#    to the gap between "session (" and ")" -- do this BEFORE the text
#    replacements below so the two runs around it stay split apart
#    instead of being merged back together by the replace.
# ---------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("session \(", $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
$anchor.Collapse(0)
$d.Bookmarks.Add("_GoBack", $anchor) | Out-Null

# ---------------------------------------------------------------------
# 2) Plain text swaps (Find/Replace over the whole document).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Going in tower…", $true, $false, $false, $false, $false, $true, 1, $false, "Ventuum tower complete", 2) | Out-Null

$d.Content.Find.Execute("Last redacted 13.05.19)", $true, $false, $false, $false, $false, $true, 1, $false, ")", 2) | Out-Null

$d.Content.Find.Execute("Disrupt new church building near the Mologia-Alliance border. (lasts 3 months)", $true, $false, $false, $false, $false, $true, 1, $false, " ", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Insert a new, empty bulleted paragraph right after the
#    "Ventuum tower complete" item (same list style/numbering as its
#    neighbours).
# ---------------------------------------------------------------------
$towerPara = $d.Content.Find.Execute("Ventuum tower complete", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$towerRange = $d.Content
$towerRange.Find.Execute("Ventuum tower complete", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$towerRange.Paragraphs.Item(1).Range.InsertParagraphAfter()

# ---------------------------------------------------------------------
# 4) Delete the three trailing questline bullets entirely (text +
#    paragraph mark), working from the bottom up so paragraph indices
#    stay valid while we go.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Zingar questline is frozen for 2-3 months, because otherwise it’ll reveal itself to publicity, which is not acceptable right now for Alfrani." -or `
        $t -eq "Help Alliance’s allies (travel to the land of tieflings and dragonborn)" -or `
        $t -eq "Help halfling allies ") {
        $p.Range.Delete() | Out-Null
    }
}
